# Update Name of Algo
# Apply updated imputed values to specific cells in columns B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 5.027500000000004
$ws.Range("C4").Value  = -14.3389

$ws.Range("C5").Value  = -14.73840000000001

$ws.Range("B6").Value  = 9.0512

$ws.Range("B7").Value  = 6.324999999999997

$ws.Range("C8").Value  = -11.86849999999999

$ws.Range("B16").Value = 8.616200000000006
$ws.Range("C16").Value = -11.9672

$ws.Range("B20").Value = 5.732699999999999

$ws.Range("C22").Value = -11.01909999999999
